# semana 38 de 2025
# Adds a new "week 38" column (AO) to the IRA hospitalario weekly tracking
# sheet, mirroring the existing week columns, and corrects the split
# between week 37 (AN) and week 38 (AO) for row 35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: AO1 = "38" (a text label, like the other week-number headers) ---
$ws.Range("AO1").NumberFormat = "@"
$ws.Range("AO1").Value = "38"

# --- Row 35 had part of its week-38 count still lumped into week 37 (AN35). ---
# Split it: AN35 11 -> 2, and put the remainder (5) into the new AO35 cell.
$ws.Range("AN35").Value = 2
$ws.Range("AO35").Value = 5

# --- New week-38 values (column AO) for every other reporting facility ---
$weekValues = [ordered]@{
    2  = 0
    5  = 0
    6  = 25
    7  = 0
    8  = 19
    9  = 0
    10 = 0
    11 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    23 = 0
    25 = 0
    28 = 3
    29 = 3
    31 = 0
    36 = 0
    37 = 0
    38 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 0
    52 = 0
    53 = 0
    54 = 0
    55 = 0
    56 = 0
    57 = 0
    58 = 0
}

foreach ($row in $weekValues.Keys) {
    $ws.Range("AO$row").Value = $weekValues[$row]
}
